$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits: reworded test-step / expected-result text (row 9 and row 28) ---
# Order matters for how new shared strings are appended to the string table,
# so we touch them in the same order the new text appears in the target file.
$ws.Range("I9").Value = "User Should be add alpha numerical text "
$ws.Range("J9").Value = "User is able to enter alpha Numberical text in Text Field"

$ws.Range("J28").Value = "User is not able to see blue border and a  black border is seen"
$ws.Range("I28").Value = "User should able to view blue border as per design document"

# --- Style cleanup: L23/L26/L28 used a duplicate "left/right border + wrap"
# style that only differed from the one used by L12 by a redundant applyFill
# flag. Re-apply the same border/wrap formatting so the cell is restyled
# onto the de-duplicated style. ---
foreach ($addr in @("L23", "L26", "L28")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    $cell.WrapText = $true
}

# --- Row height adjustments ---
$ws.Rows.Item(10).RowHeight = 43.5
$ws.Rows.Item(18).RowHeight = 29
$ws.Rows.Item(20).RowHeight = 29
$ws.Rows.Item(21).RowHeight = 29
$ws.Rows.Item(25).RowHeight = 29
$ws.Rows.Item(28).RowHeight = 29

# --- Update the active selection to match the saved view state ---
$ws.Range("I32").Select() | Out-Null
